$wb = $excel.ActiveWorkbook

# Rename the "Compensation" header in the "Job Applications" sheet
# (column I) to "What excites this opportunity".
$jobApps = $wb.Worksheets.Item("Job Applications")
$jobApps.Range("I1").Value = "What excites this opportunity"

# Make "Job Applications" the active sheet/tab with I2 selected
# (it was previously "Skills Tracker", selected at E4).
$jobApps.Activate()
$jobApps.Range("I2").Select() | Out-Null
